$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- helpers -----------------------------------------------------------
# A "plain" (un-styled, General-format) cell from the existing data area
# that we can copy formats from with PasteSpecial so that newly written
# cells don't pick up stray number-formats (e.g. from the text-forcing
# apostrophe trick below).
$plainCell = $ws.Range("A2")
# A styled header cell (bold, bordered, centered) to stamp onto new
# header cells.
$headerCell = $ws.Range("A1")

function Set-TextCell($row, $col, [string]$value) {
    # Leading apostrophe forces Excel to store the value as literal text
    # (t="s"/inlineStr) rather than auto-converting numbers/dates. An
    # empty string becomes a present-but-blank text cell this way too.
    $cell = $ws.Cells.Item($row, $col)
    $cell.Value = "'" + $value
    $plainCell.Copy()
    $cell.PasteSpecial(-4122)
}

function Set-HeaderCell($row, $col, [string]$value) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.Value = "'" + $value
    $headerCell.Copy()
    $cell.PasteSpecial(-4122)
}

# --- new "Notes" column header ------------------------------------------
Set-HeaderCell 1 5 "Notes"

# --- blank Notes cells for the pre-existing rows ------------------------
Set-TextCell 2 5 ""
Set-TextCell 3 5 ""

# --- new sighting rows ---------------------------------------------------
$rows = @(
    @("kiwi",        "2", "50 hicjs road",   "2025-08-21 12:07:21", ""),
    @("kiwi",        "1", "mill creek",      "2025-08-21 13:20:05", "it was small"),
    @("kiwi",        "2", "back road",       "2025-08-21 13:23:56", ""),
    @("haast eagle", "4", "mountain range",  "2025-08-21 19:16:02", "it was quite large"),
    @("kiwi",        "4", "50 hicks road",   "2025-08-21 19:26:29", "")
)

$r = 4
foreach ($row in $rows) {
    Set-TextCell $r 1 $row[0]
    Set-TextCell $r 2 $row[1]
    Set-TextCell $r 3 $row[2]
    Set-TextCell $r 4 $row[3]
    Set-TextCell $r 5 $row[4]
    $r++
}
